# "Generate Report for handback" - append the handback-status row for the
# newly-processed file 93a8ae22-3269-411d-a850-54e02cad6daa to all three
# sheets (Overview, zh-cn, de-de), mirroring the existing rows 2/3 pattern.

$wb = $excel.ActiveWorkbook

$fileId   = "93a8ae22-3269-411d-a850-54e02cad6daa"
$mdName   = "$fileId.md"
$xlfHash  = "e4301c2811b63bceb5186dd57b3f7a9d11e60c6a"
$zhXlf    = "$fileId.$xlfHash.zh-cn.xlf"
$deXlf    = "$fileId.$xlfHash.de-de.xlf"

$inSync   = "Handed back: in sync with en-US"
$include  = "Include"

$zhHandoffDt  = "2016-01-25 13:32:42"
$zhHandbackDt = "2016-01-25 13:33:25"
$deHandoffDt  = "2016-01-25 13:32:51"
$deHandbackDt = "2016-01-25 13:33:42"

# ---------------------------------------------------------------------
# Sheet 1: "Overview" -- File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $inSync
$wsZh.Range("C4").Value = $zhXlf
$wsZh.Range("D4").Value = $zhHandoffDt
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $zhXlf
$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $include

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item(4, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item(4, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item(4, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf",
    "",
    "",
    $zhXlf
) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $inSync
$wsDe.Range("C4").Value = $deXlf
$wsDe.Range("D4").Value = $deHandoffDt
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $deXlf
$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $include

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item(4, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item(4, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item(4, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf",
    "",
    "",
    $deXlf
) | Out-Null
